$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 258, shifting existing rows 258-274 down to 259-275.
$ws.Rows.Item(258).Insert()

# Populate the newly inserted row 258 with the new record.
$ws.Range("A258").Value = 10
$ws.Range("B258").Value = "Vega Modelo de Temuco"
$ws.Range("C258").Value = "La Araucanía"
$ws.Range("D258").Value = 45021
$ws.Range("E258").Value = 9
$ws.Range("F258").Value = 100112013
$ws.Range("G258").Value = "Alcachofa"
$ws.Range("H258").Value = "Madrigal"
$ws.Range("I258").Value = "Extra"
$ws.Range("J258").Value = 65
$ws.Range("K258").Value = 30000
$ws.Range("L258").Value = 30000
$ws.Range("M258").Value = 30000
$ws.Range("N258").Value = "$/caja 35 unidades"
$ws.Range("O258").Value = "Provincia de Limarí"
$ws.Range("P258").Value = 857
$ws.Range("Q258").Value = 35
$ws.Range("R258").Value = "Hortaliza"
